# "fixed horizontal centering on registers"
#
# Turn on "center on page horizontally" for the register sheet's print
# setup, and leave the active selection on the "Преподаватели" row
# (A7:M7) instead of the OCR-code cell (M4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Center the printed sheet horizontally on the page.
$ws.PageSetup.CenterHorizontally = $true

# Move the active selection to A7:M7 (the "Преподаватели:" merged row).
$ws.Range("A7:M7").Select()
